$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with the latest
# crypto snapshot, as produced by the scheduled GitHub Actions run.
# Values are written with a leading-quote so Excel stores them as text
# (matching the sheet's existing inline-string cells) instead of
# re-interpreting number-shaped strings like "331.88" as numerics.
$ws.Range("D2").Value = "`'30.130.63"
$ws.Range("E2").Value = "`'  +5.79%  "
$ws.Range("D3").Value = "`'1.924.61"
$ws.Range("E3").Value = "`'  +2.81%  "
$ws.Range("D5").Value = "`'331.88"
$ws.Range("E5").Value = "`'  +5.17%  "
$ws.Range("D6").Value = "`'0.9997"
$ws.Range("E6").Value = "`'  -0.77%  "
$ws.Range("D7").Value = "`'0.5241"
$ws.Range("E7").Value = "`'  +3.02%  "
$ws.Range("D8").Value = "`'0.4095"
$ws.Range("E8").Value = "`'  +5.12%  "
$ws.Range("D9").Value = "`'0.08555"
$ws.Range("E9").Value = "`'  +2.59%  "
$ws.Range("D10").Value = "`'43.54"
$ws.Range("E10").Value = "`'  +4.36%  "
$ws.Range("D11").Value = "`'1.132"
$ws.Range("E11").Value = "`'  +2.72%  "
$ws.Range("D12").Value = "`'22.53"
$ws.Range("E12").Value = "`'  +10.81%  "
$ws.Range("D13").Value = "`'6.445"
$ws.Range("E13").Value = "`'  +3.68%  "
$ws.Range("D14").Value = "`'1.920.61"
$ws.Range("E14").Value = "`'  +2.78%  "
$ws.Range("D15").Value = "`'7.425"
$ws.Range("E15").Value = "`'  +2.19%  "
$ws.Range("E16").Value = "`'  -0.78%  "
$ws.Range("D17").Value = "`'96.69"
$ws.Range("E17").Value = "`'  +6.12%  "
$ws.Range("D18").Value = "`'0.00001119"
$ws.Range("E18").Value = "`'  +1.34%  "
$ws.Range("D19").Value = "`'0.06716"
$ws.Range("E19").Value = "`'  -0.08%  "
$ws.Range("D20").Value = "`'18.38"
$ws.Range("E20").Value = "`'  +3.88%  "
$ws.Range("D21").Value = "`'1.000"
$ws.Range("E22").Value = "`'  +2.86%  "
$ws.Range("D23").Value = "`'30.136.74"
$ws.Range("E23").Value = "`'  +5.70%  "
$ws.Range("D24").Value = "`'11.31"
$ws.Range("E24").Value = "`'  +1.81%  "
$ws.Range("D25").Value = "`'2.215"
$ws.Range("E25").Value = "`'  -0.55%  "
$ws.Range("D26").Value = "`'2.141.91"
$ws.Range("E26").Value = "`'  +2.84%  "
$ws.Range("D27").Value = "`'21.24"
$ws.Range("E27").Value = "`'  +3.01%  "
$ws.Range("D28").Value = "`'160.03"
$ws.Range("E28").Value = "`'  -0.99%  "
$ws.Range("D29").Value = "`'2.482"
$ws.Range("E29").Value = "`'  +3.02%  "
$ws.Range("D30").Value = "`'130.25"
$ws.Range("E30").Value = "`'  +3.26%  "
$ws.Range("D31").Value = "`'1.086"
$ws.Range("E31").Value = "`'  +4.98%  "
$ws.Range("D32").Value = "`'0.1059"
$ws.Range("E32").Value = "`'  +1.71%  "
$ws.Range("D33").Value = "`'6.131"
$ws.Range("E33").Value = "`'  +6.92%  "
$ws.Range("D34").Value = "`'3.642"
$ws.Range("E34").Value = "`'  +1.19%  "
$ws.Range("D35").Value = "`'0.02520"
$ws.Range("E35").Value = "`'  +2.67%  "
$ws.Range("D36").Value = "`'0.06616"
$ws.Range("E36").Value = "`'  +1.16%  "
$ws.Range("D37").Value = "`'0.2227"
$ws.Range("E37").Value = "`'  +3.37%  "
$ws.Range("D38").Value = "`'1.244"
$ws.Range("E38").Value = "`'  +4.94%  "
$ws.Range("D39").Value = "`'9.051"
$ws.Range("E39").Value = "`'  +2.68%  "
$ws.Range("D40").Value = "`'5.223"
$ws.Range("E40").Value = "`'  +4.12%  "
$ws.Range("D41").Value = "`'0.6568"
$ws.Range("E41").Value = "`'  +3.22%  "
$ws.Range("E42").Value = "`'  +5.72%  "
$ws.Range("D43").Value = "`'1.240"
$ws.Range("E43").Value = "`'  -0.18%  "
$ws.Range("D44").Value = "`'0.6201"
$ws.Range("E44").Value = "`'  +3.57%  "
$ws.Range("D45").Value = "`'13.39"
$ws.Range("E45").Value = "`'  +2.79%  "
$ws.Range("D46").Value = "`'3.783"
$ws.Range("E46").Value = "`'  +2.59%  "
$ws.Range("D47").Value = "`'2.098"
$ws.Range("E47").Value = "`'  +4.89%  "
$ws.Range("E48").Value = "`'  +2.98%  "
$ws.Range("D49").Value = "`'125.24"
$ws.Range("E50").Value = "`'  +1.46%  "
$ws.Range("D51").Value = "`'80.26"
$ws.Range("E51").Value = "`'  +5.49%  "

# Re-normalize the touched cells back to the default "Normal" style so the
# quote-prefix/text formatting used above to force text doesn't linger on
# the cells themselves.
$ws.Range("D2","E2","D3","E3","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E22","D23","E23","D24","E24","D25","E25","D26","E26","D27","E27","D28","E28","D29","E29","D30","E30","D31","E31","D32","E32","D33","E33","D34","E34","D35","E35","D36","E36","D37","E37","D38","E38","D39","E39","D40","E40","D41","E41","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","E48","D49","E50","D51","E51").Style = "Normal"
